$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.288577556610107
$ws.Range("B1").Value = 2.334793567657471
$ws.Range("C1").Value = 6.292636871337891
$ws.Range("D1").Value = 1.575372099876404
$ws.Range("E1").Value = 1.331328630447388
